# WTREGEN.xlsx weekly refresh
#
# Appends the two newest weekly observations to the "Data" sheet and syncs
# the "SeriesInfo" metadata sheet (realtime_start/realtime_end,
# observation_end, last_updated) to match the refreshed FRED pull.

$wb = $excel.ActiveWorkbook

$dataSheet  = $wb.Worksheets.Item("Data")
$infoSheet  = $wb.Worksheets.Item("SeriesInfo")

# --- Data sheet: append rows 98 and 99, inheriting the date formatting ---
# used by the rest of column A (copy format from the last existing row,
# then overwrite with the new date-serial / value pair).

$dataSheet.Range("A97").Copy()
$dataSheet.Range("A98").PasteSpecial(-4122)
$dataSheet.Range("A98").Value = 45147
$dataSheet.Range("B98").Value = 444.594

$dataSheet.Range("A97").Copy()
$dataSheet.Range("A99").PasteSpecial(-4122)
$dataSheet.Range("A99").Value = 45154
$dataSheet.Range("B99").Value = 435.795

# --- SeriesInfo sheet: refresh metadata rows ---
# realtime_start / realtime_end / observation_end are plain "yyyy-mm-dd"
# text; Excel's smart entry otherwise re-interprets that text as a date
# serial, so force the literal text for those three cells.
$infoSheet.Range("B3").Formula = "'2023-08-22"
$infoSheet.Range("B4").Formula = "'2023-08-22"
$infoSheet.Range("B7").Formula = "'2023-08-16"

# last_updated already carries a time + UTC-offset suffix, which Excel's
# parser does not treat as a recognizable date/time literal, so a plain
# value assignment keeps it as text.
$infoSheet.Range("B14").Value = "2023-08-17 15:35:18-05"
